$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the Coin/Link/Price/Volume(1h) table with the latest feed values.
$ws.Range("D2").Value = "'26.583.36"
$ws.Range("E2").Value = "  +0.73%  "
$ws.Range("D3").Value = "'1.820.93"
$ws.Range("E3").Value = "  +1.48%  "
$ws.Range("D4").Value = "'1.010"
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "'1.008"
$ws.Range("E5").Value = "  +0.22%  "
$ws.Range("D6").Value = "'305.01"
$ws.Range("E6").Value = "  -0.68%  "
$ws.Range("D7").Value = "'0.4642"
$ws.Range("D8").Value = "'0.3591"
$ws.Range("E8").Value = "  -0.22%  "
$ws.Range("D9").Value = "'0.07124"
$ws.Range("E9").Value = "  +0.09%  "
$ws.Range("D10").Value = "'0.8967"
$ws.Range("E10").Value = "  +1.46%  "
$ws.Range("D11").Value = "'0.07772"
$ws.Range("E11").Value = "  -0.61%  "
$ws.Range("D12").Value = "'19.27"
$ws.Range("E12").Value = "  -1.15%  "
$ws.Range("D13").Value = "'1.815.11"
$ws.Range("E13").Value = "  +0.90%  "
$ws.Range("D14").Value = "'5.239"
$ws.Range("E14").Value = "  -0.85%  "
$ws.Range("D15").Value = "'6.323"
$ws.Range("E15").Value = "  -0.02%  "
$ws.Range("D16").Value = "'87.22"
$ws.Range("E16").Value = "  +2.52%  "
$ws.Range("D17").Value = "'1.011"
$ws.Range("E17").Value = "  +0.22%  "
$ws.Range("D18").Value = "'0.000008529"
$ws.Range("E18").Value = "  -0.58%  "
$ws.Range("E19").Value = "  +0.24%  "
$ws.Range("D20").Value = "'26.626.87"
$ws.Range("E20").Value = "  +0.76%  "
$ws.Range("D21").Value = "'14.12"
$ws.Range("E21").Value = "  -1.05%  "
$ws.Range("D22").Value = "'5.005"
$ws.Range("E22").Value = "  +0.33%  "
$ws.Range("D23").Value = "'10.51"
$ws.Range("E23").Value = "  -0.20%  "
$ws.Range("D24").Value = "'1.914"
$ws.Range("E24").Value = "  -3.15%  "
$ws.Range("D25").Value = "'151.98"
$ws.Range("E25").Value = "  -0.35%  "
$ws.Range("D26").Value = "'17.87"
$ws.Range("E26").Value = "  -0.15%  "
$ws.Range("D27").Value = "'1.963"
$ws.Range("E27").Value = "  -4.11%  "
$ws.Range("D28").Value = "'113.50"
$ws.Range("E28").Value = "  +1.42%  "
$ws.Range("D29").Value = "'4.789"
$ws.Range("E29").Value = "  -1.70%  "
$ws.Range("D30").Value = "'0.08792"
$ws.Range("E30").Value = "  +1.58%  "
$ws.Range("D31").Value = "'3.129"
$ws.Range("E31").Value = "  +2.44%  "
$ws.Range("D32").Value = "'0.7266"
$ws.Range("E32").Value = "  +0.04%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "'4.423"
$ws.Range("E33").Value = "  -0.59%  "
$ws.Range("B34").Value = "RenderToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D34").Value = "'2.716"
$ws.Range("E34").Value = "  +0.20%  "
$ws.Range("D35").Value = "'1.121"
$ws.Range("E35").Value = "  +1.03%  "
$ws.Range("D36").Value = "'1.071"
$ws.Range("E36").Value = "  -0.40%  "
$ws.Range("D37").Value = "'0.01918"
$ws.Range("E37").Value = "  -1.28%  "
$ws.Range("E38").Value = "  +1.44%  "
$ws.Range("D39").Value = "'0.05099"
$ws.Range("E39").Value = "  -0.24%  "
$ws.Range("D40").Value = "'6.834"
$ws.Range("E40").Value = "  -0.91%  "
$ws.Range("D41").Value = "'0.5029"
$ws.Range("E41").Value = "  -2.82%  "
$ws.Range("D42").Value = "'0.1489"
$ws.Range("E42").Value = "  -2.40%  "
$ws.Range("E43").Value = "  -0.84%  "
$ws.Range("D44").Value = "'1.009"
$ws.Range("E44").Value = "  +0.24%  "
$ws.Range("D45").Value = "'0.4636"
$ws.Range("E45").Value = "  -0.75%  "
$ws.Range("D46").Value = "'9.932"
$ws.Range("E46").Value = "  +0.65%  "
$ws.Range("D47").Value = "'97.90"
$ws.Range("E47").Value = "  -2.27%  "
$ws.Range("D48").Value = "'1.552"
$ws.Range("E48").Value = "  -2.24%  "
$ws.Range("D49").Value = "'0.05983"
$ws.Range("E49").Value = "  +0.14%  "
$ws.Range("D50").Value = "'63.58"
$ws.Range("E50").Value = "  -1.00%  "
$ws.Range("D51").Value = "'35.66"
$ws.Range("E51").Value = "  -1.76%  "

# Clear the quote-prefix style the apostrophe above introduced so the
# cell format matches the rest of the (unstyled) data cells.
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
